$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Datos actualizados" timestamp in A1 (07:40 -> 08:57)
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 08:57"

# 2. Update Ucrania row (row 27) statistics
$ws.Range("B27").Value = 140479
$ws.Range("C27").Value = 2411
$ws.Range("D27").Value = 63546
$ws.Range("E27").Value = 73999
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 57
$ws.Range("H27").Value = 2934

# 3. Swap Uruguay / Georgia: Georgia now appears before Uruguay (row 151),
#    with freshly updated figures, while Uruguay (row 152) keeps the
#    previous Uruguay figures that used to live in row 151.
$ws.Range("A151").Value = "Georgia"
$ws.Range("B151").Value = 1729
$ws.Range("C151").Value = 45
$ws.Range("D151").Value = 1321
$ws.Range("E151").Value = 389
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 19

$ws.Range("A152").Value = "Uruguay"
$ws.Range("B152").Value = 1693
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 1466
$ws.Range("E152").Value = 182
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 45
